$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "New England Shelf" cruise label to "North East Shelf"
# throughout column E (Cruise).
$range = $ws.Range("E1:E97")
$range.Replace("New England Shelf", "North East Shelf")

# Update the active cell selection left behind by the editing session.
$ws.Range("L6").Select()
